$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 838.38464
$ws.Range("I32").Value = 600
$ws.Range("K32").Value = 600
$ws.Range("M32").Value = -274

$ws.Range("H39").Value = 216.3125
$ws.Range("I39").Value = 164.06667
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 492.20001
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -196.20001
$ws.Range("N39").Value = -3592

$ws.Range("H51").Value = 5141.5713
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5141.5713
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5141.5713
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -6109.5713

$ws.Range("H53").Value = 329.56
$ws.Range("J53").Value = 169.15384
$ws.Range("L53").Value = 169.15384
$ws.Range("N53").Value = -1443.15384

$ws.Range("H98").Value = 2484.8462
$ws.Range("I98").Value = 1394.5
$ws.Range("K98").Value = 1394.5
$ws.Range("M98").Value = 103.5

$ws.Range("H122").Value = 2484.8462
$ws.Range("I122").Value = 1394.5
$ws.Range("K122").Value = 4183.5
$ws.Range("M122").Value = -1733.5

$ws.Range("H141").Value = 5661.727
$ws.Range("I141").Value = 5661.727
$ws.Range("K141").Value = 16985.181
$ws.Range("M141").Value = -11805.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996

$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984

$ws.Range("H92").Value = 39275
$ws.Range("J92").Value = 39275
$ws.Range("L92").Value = 39275
$ws.Range("N92").Value = -44267

$ws.Range("H122").Value = 4214
$ws.Range("I122").Value = 4583.1665
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 13749.4995
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -11299.4995
$ws.Range("N122").Value = -10897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 19542.5
$ws.Range("J76").Value = 19542.5
$ws.Range("L76").Value = 19542.5
$ws.Range("N76").Value = -20172.5

$ws.Range("H79").Value = 19542.5
$ws.Range("J79").Value = 19542.5
$ws.Range("L79").Value = 19542.5
$ws.Range("N79").Value = -21726.5

$ws.Range("H94").Value = 2782.1667
$ws.Range("I94").Value = 2673.75
$ws.Range("J94").Value = 2999
$ws.Range("K94").Value = 2673.75
$ws.Range("L94").Value = 2999
$ws.Range("M94").Value = -2222.75
$ws.Range("N94").Value = -3901

$ws.Range("H95").Value = 34000
$ws.Range("J95").Value = 34000
$ws.Range("L95").Value = 34000
$ws.Range("N95").Value = -39492

$ws.Range("H103").Value = 19966.666
$ws.Range("J103").Value = 19966.666
$ws.Range("L103").Value = 19966.666
$ws.Range("N103").Value = -22310.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 13887.444
$ws.Range("J88").Value = 13887.444
$ws.Range("L88").Value = 13887.444
$ws.Range("N88").Value = -14699.444

$ws.Range("H91").Value = 13887.444
$ws.Range("J91").Value = 13887.444
$ws.Range("L91").Value = 13887.444
$ws.Range("N91").Value = -16695.444

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 981.6
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 477
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 1431
$ws.Range("M97").Value = -8504
$ws.Range("N97").Value = -2423

$ws.Range("H113").Value = 2039.8
$ws.Range("J113").Value = 2039.8
$ws.Range("L113").Value = 6119.4
$ws.Range("N113").Value = -10459.4

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 58333.332
$ws.Range("J34").Value = 58333.332
$ws.Range("L34").Value = 58333.332
$ws.Range("N34").Value = -58869.332

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H76").Value = 58333.332
$ws.Range("J76").Value = 58333.332
$ws.Range("L76").Value = 58333.332
$ws.Range("N76").Value = -58963.332

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H79").Value = 58333.332
$ws.Range("J79").Value = 58333.332
$ws.Range("L79").Value = 58333.332
$ws.Range("N79").Value = -60517.332

$ws.Range("H80").Value = 8277.177
$ws.Range("I80").Value = 7749.778
$ws.Range("J80").Value = 8870.5
$ws.Range("K80").Value = 7749.778
$ws.Range("L80").Value = 8870.5
$ws.Range("M80").Value = -6751.778
$ws.Range("N80").Value = -10866.5

$ws.Range("H83").Value = 8277.177
$ws.Range("I83").Value = 7749.778
$ws.Range("J83").Value = 8870.5
$ws.Range("K83").Value = 38748.89
$ws.Range("L83").Value = 44352.5
$ws.Range("M83").Value = -33756.89
$ws.Range("N83").Value = -54336.5

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H126").Value = 4598.3335
$ws.Range("I126").Value = 4597.5
$ws.Range("K126").Value = 13792.5
$ws.Range("M126").Value = -11322.5

$ws.Range("H130").Value = 78000
$ws.Range("J130").Value = 78000
$ws.Range("L130").Value = 78000
$ws.Range("N130").Value = -88040

$ws.Range("H132").Value = 2225.2
$ws.Range("I132").Value = 2219.5
$ws.Range("J132").Value = 2248
$ws.Range("K132").Value = 6658.5
$ws.Range("L132").Value = 6744
$ws.Range("M132").Value = -4128.5
$ws.Range("N132").Value = -11804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H40").Value = 5370.2856
$ws.Range("I40").Value = 5448.8335
$ws.Range("K40").Value = 5448.8335
$ws.Range("M40").Value = -5312.8335

$ws.Range("H46").Value = 2611
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2745.25
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2745.25
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -3121.25

$ws.Range("H68").Value = 599.6667
$ws.Range("I68").Value = 599.6667
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 599.6667
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 149.3333
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 599.6667
$ws.Range("I71").Value = 599.6667
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 2998.3335
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 745.6665000000003
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 2958.818
$ws.Range("I82").Value = 2309.4
$ws.Range("K82").Value = 2309.4
$ws.Range("M82").Value = -1948.4

$ws.Range("H85").Value = 2958.818
$ws.Range("I85").Value = 2309.4
$ws.Range("K85").Value = 2309.4
$ws.Range("M85").Value = -1061.4

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50856

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52964

$ws.Range("H110").Value = 49999
$ws.Range("J110").Value = 49999
$ws.Range("L110").Value = 49999
$ws.Range("N110").Value = -58179

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 224999
$ws.Range("I70").Value = 200000
$ws.Range("J70").Value = 249998
$ws.Range("K70").Value = 200000
$ws.Range("L70").Value = 249998
$ws.Range("M70").Value = -199685
$ws.Range("N70").Value = -250628

$ws.Range("H73").Value = 224999
$ws.Range("I73").Value = 200000
$ws.Range("J73").Value = 249998
$ws.Range("K73").Value = 200000
$ws.Range("L73").Value = 249998
$ws.Range("M73").Value = -198908
$ws.Range("N73").Value = -252182

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H107").Value = 412.7143
$ws.Range("I107").Value = 422.5
$ws.Range("J107").Value = 399.66666
$ws.Range("K107").Value = 1267.5
$ws.Range("L107").Value = 1198.99998
$ws.Range("M107").Value = 652.5
$ws.Range("N107").Value = -5038.999980000001

$ws.Range("H122").Value = 2874.5
$ws.Range("I122").Value = 1499
$ws.Range("K122").Value = 4497
$ws.Range("M122").Value = -2047
